$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Friday (column H) as "done" for week row 21, matching the green
# fill/border formatting already used for other completed day cells
# (e.g. H19/H20), and bump that week's total day count by one.
$ws.Range("H20").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("I21").Value = 5

# Recalculate so the dependent formulas (K3, L3, M3) refresh their cached
# values to reflect the new total.
$excel.Calculate()

# Move the active selection from L20 to A2.
$ws.Range("A2").Select() | Out-Null
